$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1835205992509363
$ws.Range("C2").Value = 0.5730337078651685
$ws.Range("J2").Value = 0.003745318352059925
$ws.Range("P2").Value = 0.1423220973782771
$ws.Range("S2").Value = 0.09737827715355805
$ws.Range("C3").Value = 0.03773584905660377
$ws.Range("J3").Value = 0.02515723270440252
$ws.Range("P3").Value = 0.7358490566037735
$ws.Range("S3").Value = 0.2012578616352201
$ws.Range("J4").Value = 0.05714285714285714
$ws.Range("P4").Value = 0.7428571428571429
$ws.Range("S4").Value = 0.2
$ws.Range("B6").Value = 0.06926406926406926
$ws.Range("F6").Value = 0.06060606060606061
$ws.Range("J6").Value = 0.2424242424242424
$ws.Range("O6").Value = 0.03463203463203463
$ws.Range("Q6").Value = 0.1168831168831169
$ws.Range("R6").Value = 0.05194805194805195
$ws.Range("S6").Value = 0.4242424242424243
$ws.Range("B7").Value = 0.08968609865470852
$ws.Range("D7").Value = 0.0179372197309417
$ws.Range("F7").Value = 0.06278026905829596
$ws.Range("J7").Value = 0.1479820627802691
$ws.Range("O7").Value = 0.008968609865470852
$ws.Range("Q7").Value = 0.1345291479820628
$ws.Range("R7").Value = 0.05829596412556054
$ws.Range("S7").Value = 0.4798206278026906
$ws.Range("B8").Value = 0.08656036446469248
$ws.Range("D8").Value = 0.009111617312072893
$ws.Range("E8").Value = 0.002277904328018223
$ws.Range("F8").Value = 0.05922551252847381
$ws.Range("J8").Value = 0.132118451025057
$ws.Range("O8").Value = 0.02050113895216401
$ws.Range("Q8").Value = 0.1526195899772209
$ws.Range("R8").Value = 0.07289293849658314
$ws.Range("S8").Value = 0.4646924829157175
$ws.Range("B9").Value = 0.07106598984771574
$ws.Range("D9").Value = 0.01522842639593909
$ws.Range("F9").Value = 0.05076142131979695
$ws.Range("J9").Value = 0.07614213197969544
$ws.Range("O9").Value = 0.01015228426395939
$ws.Range("Q9").Value = 0.16751269035533
$ws.Range("R9").Value = 0.1065989847715736
$ws.Range("S9").Value = 0.5025380710659898
$ws.Range("B10").Value = 0.1002444987775061
$ws.Range("D10").Value = 0.01792991035044825
$ws.Range("F10").Value = 0.07090464547677261
$ws.Range("J10").Value = 0.1157294213528932
$ws.Range("O10").Value = 0.01385493072534637
$ws.Range("Q10").Value = 0.1792991035044825
$ws.Range("R10").Value = 0.05215973920130399
$ws.Range("S10").Value = 0.449877750611247
$ws.Range("G11").Value = 0.1614730878186969
$ws.Range("J11").Value = 0.0764872521246459
$ws.Range("K11").Value = 0.2011331444759207
$ws.Range("L11").Value = 0.5382436260623229
$ws.Range("S11").Value = 0.0226628895184136
$ws.Range("G12").Value = 0.7760416666666666
$ws.Range("J12").Value = 0.1666666666666667
$ws.Range("K12").Value = 0.005208333333333333
$ws.Range("L12").Value = 0.015625
$ws.Range("S12").Value = 0.03645833333333334
$ws.Range("G13").Value = 0.7567567567567568
$ws.Range("J13").Value = 0.1891891891891892
$ws.Range("S13").Value = 0.05405405405405406
$ws.Range("F15").Value = 0.02403846153846154
$ws.Range("H15").Value = 0.1394230769230769
$ws.Range("I15").Value = 0.1057692307692308
$ws.Range("J15").Value = 0.3317307692307692
$ws.Range("K15").Value = 0.07211538461538461
$ws.Range("M15").Value = 0.01442307692307692
$ws.Range("O15").Value = 0.05288461538461538
$ws.Range("S15").Value = 0.2596153846153846
$ws.Range("F16").Value = 0.01666666666666667
$ws.Range("H16").Value = 0.1277777777777778
$ws.Range("I16").Value = 0.09444444444444444
$ws.Range("J16").Value = 0.4333333333333333
$ws.Range("K16").Value = 0.1
$ws.Range("O16").Value = 0.07222222222222222
$ws.Range("S16").Value = 0.1555555555555556
$ws.Range("F17").Value = 0.02972972972972973
$ws.Range("H17").Value = 0.1540540540540541
$ws.Range("I17").Value = 0.0918918918918919
$ws.Range("J17").Value = 0.4081081081081081
$ws.Range("K17").Value = 0.1054054054054054
$ws.Range("M17").Value = 0.01351351351351351
$ws.Range("O17").Value = 0.05945945945945946
$ws.Range("S17").Value = 0.1378378378378378
$ws.Range("F18").Value = 0.02836879432624113
$ws.Range("H18").Value = 0.2127659574468085
$ws.Range("I18").Value = 0.07801418439716312
$ws.Range("J18").Value = 0.3900709219858156
$ws.Range("K18").Value = 0.1205673758865248
$ws.Range("M18").Value = 0.02127659574468085
$ws.Range("O18").Value = 0.02836879432624113
$ws.Range("S18").Value = 0.1205673758865248
$ws.Range("F19").Value = 0.01939799331103679
$ws.Range("H19").Value = 0.2053511705685619
$ws.Range("I19").Value = 0.07357859531772576
$ws.Range("J19").Value = 0.348494983277592
$ws.Range("K19").Value = 0.1224080267558528
$ws.Range("M19").Value = 0.01672240802675585
$ws.Range("N19").Value = 0.001337792642140468
$ws.Range("O19").Value = 0.06421404682274247
$ws.Range("S19").Value = 0.148494983277592
